$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.238453984260559
$ws.Range("B1").Value = 2.368208169937134
$ws.Range("C1").Value = 3.236217498779297
$ws.Range("D1").Value = 3.527665853500366
$ws.Range("E1").Value = 1.10883367061615
